# Updated cryptos list on Wed Jul  3 11:36:48 UTC 2024 with GitHub Actions
#
# Refreshes Price (D) / Volume(1h) (E) columns with newly scraped values for
# most rows; rows 37-38 also swap which coin (Monero/Aptos) occupies each
# ranking slot, so Coin (B) and Link (C) are updated there too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E store plain-looking numbers as TEXT in this sheet (e.g. "1.00",
# "0.119", "  -3.81%  ") so formatting/precision survives re-export. Force
# text format first on any new value that Excel would otherwise silently
# reinterpret as a number (and thus mangle via float rounding).

$ws.Range("D2").Value = "60.295.01"
$ws.Range("E2").Value = "  -3.81%  "
$ws.Range("D3").Value = "3.303.34"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.28"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.68"
$ws.Range("E6").Value = "  -3.55%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.310.53"
$ws.Range("E8").Value = "  -3.89%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.119"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.407"
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("D13").Value = "3.871.92"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.21"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "3.316.93"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "60.355.61"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.55"
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.81"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.91"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.548"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -8.72%  "
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  -6.31%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.61"
$ws.Range("E33").Value = "  -5.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.55"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "165.93"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.74"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("E39").Value = "  -6.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.51"
$ws.Range("E40").Value = "  -14.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0737"
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.95"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.752"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.20"
$ws.Range("E45").Value = "  -4.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.59"
$ws.Range("E46").Value = "  -5.31%  "
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D48").Value = "2.372.37"
$ws.Range("E48").Value = "  -7.65%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.58"
$ws.Range("E50").Value = "  -4.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.62"
$ws.Range("E51").Value = "  -4.18%  "
